$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "ID"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "LName "
$ws.Range("E1").Value = "Add"
$ws.Range("F1").Value = "Marks "

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "D1"
$ws.Range("C2").Value = "om"
$ws.Range("D2").Value = "na"
$ws.Range("E2").Value = "dsad"
$ws.Range("F2").Value = 89

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "D2"
$ws.Range("C3").Value = "omi"
$ws.Range("D3").Value = "nan"
$ws.Range("E3").Value = "ada"
$ws.Range("F3").Value = 89

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "D3"
$ws.Range("C4").Value = "pj"
$ws.Range("D4").Value = "annana"
$ws.Range("E4").Value = "as"
$ws.Range("F4").Value = 22

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "D4"
$ws.Range("C5").Value = "pj1"
$ws.Range("D5").Value = "nana"
$ws.Range("E5").Value = "asd"
$ws.Range("F5").Value = 28

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "D5"
$ws.Range("C6").Value = "pj11"
$ws.Range("D6").Value = "d"
$ws.Range("E6").Value = "asd"
$ws.Range("F6").Value = 90

$ws.Range("I10").Select()
